$d = $word.ActiveDocument

# "Versi" + "on" -> merge into a single run "Version" (text unchanged,
# but Find/Replace across the run boundary coalesces the split run).
$d.Content.Find.Execute("on", $false, $false, $false, $false, $false, $true, 1, $false, "on", 2)

# " 2" -> " 1." ; stop before the "_GoBack" bookmark so it is preserved
# in place (between the "1." run and the old trailing "." run).
$d.Content.Find.Execute(" 2", $false, $false, $false, $false, $false, $true, 1, $false, " 1.", 2)

# The trailing "." run (now duplicated, after the bookmark) is removed.
$full = $d.Content
$text = $full.Text
$dotIndex = $text.IndexOf("..")
$r = $d.Range($dotIndex + 1, $dotIndex + 2)
$r.Delete()
